$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("T2").Value = 1.86
$ws.Range("V2").Value = 1.81
$ws.Range("F3").Value = 1.88
$ws.Range("G3").Value = 1.98
$ws.Range("J3").Value = 3.65
$ws.Range("Q3").Value = 1.97
$ws.Range("U3").Value = 1.96
$ws.Range("W3").Value = 2.02
$ws.Range("AC3").Value = 1000
$ws.Range("F4").Value = 2.46
$ws.Range("H4").Value = 2.72
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.78
$ws.Range("N4").Value = 2.64
$ws.Range("O4").Value = 1.29
$ws.Range("R4").Value = 1.21
$ws.Range("V4").Value = 1.37
$ws.Range("W4").Value = 1.44
$ws.Range("F5").Value = 1.5
$ws.Range("G5").Value = 1.58
$ws.Range("I5").Value = 8.4
$ws.Range("J5").Value = 4.2
$ws.Range("L5").Value = 1.32
$ws.Range("N5").Value = 3.7
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 1.87
$ws.Range("T5").Value = 1.98
$ws.Range("W5").Value = 2.7
$ws.Range("H6").Value = 2.32
$ws.Range("I6").Value = 2.52
$ws.Range("Q6").Value = 1.96
$ws.Range("U6").Value = 2.26
$ws.Range("V6").Value = 1.66
$ws.Range("W6").Value = 1.41
$ws.Range("F7").Value = 1.41
$ws.Range("G7").Value = 1.46
$ws.Range("K7").Value = 5.7
$ws.Range("N7").Value = 4.6
$ws.Range("P7").Value = 2.28
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = 2.66
$ws.Range("T7").Value = 1.91
$ws.Range("W7").Value = 3.15
$ws.Range("AA7").Value = 340
$ws.Range("AB7").Value = 9.6
$ws.Range("F8").Value = 1.63
$ws.Range("H8").Value = 6.6
$ws.Range("K8").Value = 3.85
$ws.Range("N8").Value = 2.58
$ws.Range("Q8").Value = 2.4
$ws.Range("F9").Value = 1.97
$ws.Range("G9").Value = 2.34
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 4.4
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 1.27
$ws.Range("O9").Value = 1.27
$ws.Range("R9").Value = 1.36
$ws.Range("S9").Value = 2.8
$ws.Range("T9").Value = 1.66
$ws.Range("U9").Value = 2.08
$ws.Range("V9").Value = 1.3
$ws.Range("W9").Value = 1.74
